$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return 0
}

# ---------------------------------------------------------------------
# Change 1: Title "spatial" -> "Spatial" (capitalize first letter)
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1).Range
$null = $titlePara.Find.Execute('spatial', $true, $false, $false, $false, $false, $true, 1, $false, 'Spatial', 2)

# ---------------------------------------------------------------------
# Change 2: The paragraph that holds the grid-plot picture also carries
# a trailing " " run and the "Now we have set up ... datasets." run.
# Those two runs move out into their own new Body Text paragraph placed
# right after the picture paragraph (and "reconds" is corrected to
# "records" along the way).
# ---------------------------------------------------------------------
$oldSentence = 'Now we have set up our points and a grid to interpolate onto, we are ready to do an interpolation. The method I am using is the Inverse Distance Weighting (IDW). This method is reasonable and fairly simple to use as the IDW method does not require a special setup to model spatial relationships. I use this method frequently when mapping dozens of pressure gauge locations onto a topological mapping of terrain. I am routinely doing this to process dozens of time based pressure reconds consisting of hundreds of thousands of time-value data pairs. R is a great tool to have to work with a lot of large datasets.'
$newSentence = 'Now we have set up our points and a grid to interpolate onto, we are ready to do an interpolation. The method I am using is the Inverse Distance Weighting (IDW). This method is reasonable and fairly simple to use as the IDW method does not require a special setup to model spatial relationships. I use this method frequently when mapping dozens of pressure gauge locations onto a topological mapping of terrain. I am routinely doing this to process dozens of time based pressure records consisting of hundreds of thousands of time-value data pairs. R is a great tool to have to work with a lot of large datasets.'

$pictureParaIndex = Find-ParagraphIndex $d 'Now we have set up our points'

# Strip the trailing sentence (and its leading space run) from the
# picture paragraph, leaving just the drawing behind.
$picturePara = $d.Paragraphs.Item($pictureParaIndex)
$null = $picturePara.Range.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, '', 2)
$picturePara = $d.Paragraphs.Item($pictureParaIndex)
$null = $picturePara.Range.Find.Execute(' ', $true, $false, $false, $false, $false, $true, 1, $false, '', 2)

# Insert a fresh Body Text paragraph right after the picture paragraph and
# give it the corrected sentence.
$picturePara = $d.Paragraphs.Item($pictureParaIndex)
$picturePara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($pictureParaIndex + 1)
$newPara1.Range.Style = 'Body Text'
$newPara1.Range.Text = $newSentence

# ---------------------------------------------------------------------
# Change 3: Add a brand-new Body Text paragraph right after the
# "Create the object window ... spatial.data$value." paragraph.
# ---------------------------------------------------------------------
$createParaIndex = Find-ParagraphIndex $d 'Create the object window which is defined'
$createPara = $d.Paragraphs.Item($createParaIndex)
$createPara.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($createParaIndex + 1)
$newPara2.Range.Style = 'Body Text'
$newPara2.Range.Text = 'Show the constructed grid plot with data points again for reference against the interpolated plot, which follows.'
